$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.216.53"
$ws.Range("E2").Value = "  +8.98%  "

# Row 3
$ws.Range("D3").Value = "1.758.77"
$ws.Range("E3").Value = "  +4.97%  "

# Row 4
$ws.Range("D4").Value = "'0.9958"
$ws.Range("E4").Value = "  -0.51%  "

# Row 5
$ws.Range("D5").Value = "'336.37"
$ws.Range("E5").Value = "  +1.25%  "

# Row 6
$ws.Range("D6").Value = "'0.9961"
$ws.Range("E6").Value = "  -0.24%  "

# Row 7
$ws.Range("D7").Value = "'0.3753"

# Row 8
$ws.Range("D8").Value = "'48.90"
$ws.Range("E8").Value = "  +3.67%  "

# Row 9
$ws.Range("D9").Value = "'0.3422"
$ws.Range("E9").Value = "  +5.36%  "

# Row 10
$ws.Range("D10").Value = "'1.200"
$ws.Range("E10").Value = "  +4.69%  "

# Row 11
$ws.Range("D11").Value = "'0.07589"
$ws.Range("E11").Value = "  +6.26%  "

# Row 12
$ws.Range("D12").Value = "'0.9945"
$ws.Range("E12").Value = "  -0.42%  "

# Row 13
$ws.Range("D13").Value = "'6.425"
$ws.Range("E13").Value = "  +5.35%  "

# Row 14
$ws.Range("D14").Value = "'20.65"
$ws.Range("E14").Value = "  +4.70%  "

# Row 15
$ws.Range("D15").Value = "'7.111"
$ws.Range("E15").Value = "  +6.90%  "

# Row 16
$ws.Range("D16").Value = "1.755.56"
$ws.Range("E16").Value = "  +5.27%  "

# Row 17
$ws.Range("D17").Value = "'0.00001097"
$ws.Range("E17").Value = "  +4.38%  "

# Row 18
$ws.Range("D18").Value = "'0.06727"
$ws.Range("E18").Value = "  +2.63%  "

# Row 19
$ws.Range("D19").Value = "'83.48"
$ws.Range("E19").Value = "  +5.80%  "

# Row 20
$ws.Range("D20").Value = "'0.9946"
$ws.Range("E20").Value = "  -0.44%  "

# Row 21
$ws.Range("D21").Value = "'16.92"
$ws.Range("E21").Value = "  +6.42%  "

# Row 22
$ws.Range("D22").Value = "'6.274"
$ws.Range("E22").Value = "  +5.96%  "

# Row 23
$ws.Range("D23").Value = "'12.97"
$ws.Range("E23").Value = "  +0.71%  "

# Row 24
$ws.Range("D24").Value = "27.144.26"
$ws.Range("E24").Value = "  +8.74%  "

# Row 25
$ws.Range("D25").Value = "'2.442"
$ws.Range("E25").Value = "  +0.19%  "

# Row 26
$ws.Range("D26").Value = "'1.501"
$ws.Range("E26").Value = "  +26.20%  "

# Row 27
$ws.Range("D27").Value = "'2.453"
$ws.Range("E27").Value = "  +1.96%  "

# Row 28
$ws.Range("D28").Value = "'152.50"
$ws.Range("E28").Value = "  +2.79%  "

# Row 29
$ws.Range("D29").Value = "'19.79"
$ws.Range("E29").Value = "  +5.56%  "

# Row 30
$ws.Range("D30").Value = "1.953.45"
$ws.Range("E30").Value = "  +5.51%  "

# Row 31
$ws.Range("D31").Value = "'133.30"
$ws.Range("E31").Value = "  +5.72%  "

# Row 32
$ws.Range("D32").Value = "'4.115"
$ws.Range("E32").Value = "  +0.65%  "

# Row 33
$ws.Range("D33").Value = "'6.115"
$ws.Range("E33").Value = "  +5.20%  "

# Row 34
$ws.Range("D34").Value = "'0.08657"

# Row 35
$ws.Range("E35").Value = "  +5.55%  "

# Row 36
$ws.Range("D36").Value = "'1.692"
$ws.Range("E36").Value = "  +2.00%  "

# Row 37
$ws.Range("D37").Value = "'5.514"
$ws.Range("E37").Value = "  +6.60%  "

# Row 38
$ws.Range("D38").Value = "'0.02372"
$ws.Range("E38").Value = "  +5.70%  "

# Row 39
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.06361"
$ws.Range("E39").Value = "  +4.73%  "

# Row 40
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").Value = "'0.2200"
$ws.Range("E40").Value = "  +5.18%  "

# Row 41
$ws.Range("D41").Value = "'8.622"
$ws.Range("E41").Value = "  +4.28%  "

# Row 42
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "'0.6424"
$ws.Range("E42").Value = "  +7.37%  "

# Row 43
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'1.227"
$ws.Range("E43").Value = "  -0.26%  "

# Row 44
$ws.Range("D44").Value = "'14.45"
$ws.Range("E44").Value = "  +5.55%  "

# Row 45
$ws.Range("D45").Value = "'0.9956"
$ws.Range("E45").Value = "  -0.24%  "

# Row 46
$ws.Range("D46").Value = "'0.6313"
$ws.Range("E46").Value = "  +10.00%  "

# Row 47
$ws.Range("D47").Value = "'3.944"
$ws.Range("E47").Value = "  +2.34%  "

# Row 48
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'2.104"
$ws.Range("E48").Value = "  +6.91%  "

# Row 49
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "'130.50"
$ws.Range("E49").Value = "  +4.55%  "

# Row 50
$ws.Range("D50").Value = "'0.07248"
$ws.Range("E50").Value = "  +3.47%  "

# Row 51
$ws.Range("D51").Value = "'79.08"
$ws.Range("E51").Value = "  +5.90%  "
